$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1.16
$ws.Range("E3").Value = 0.68
$ws.Range("H3").Value = 1.76
$ws.Range("K3").Value = 1.13

$ws.Range("B4").Value = 1.15
$ws.Range("E4").Value = 0.7
$ws.Range("H4").Value = 1.8
$ws.Range("K4").Value = 1.17

$ws.Range("B5").Value = 1.15
$ws.Range("E5").Value = 0.73
$ws.Range("H5").Value = 1.84
$ws.Range("K5").Value = 1.22

$ws.Range("B6").Value = 1.16
$ws.Range("E6").Value = 0.76
$ws.Range("H6").Value = 1.89
$ws.Range("K6").Value = 1.27

$ws.Range("B7").Value = 1.17
$ws.Range("E7").Value = 0.8
$ws.Range("H7").Value = 1.94
$ws.Range("K7").Value = 1.32

$ws.Range("B8").Value = 1.2
$ws.Range("E8").Value = 0.85
$ws.Range("H8").Value = 2.01
$ws.Range("K8").Value = 1.38

$ws.Range("B9").Value = 1.23
$ws.Range("E9").Value = 0.9
$ws.Range("H9").Value = 2.08
$ws.Range("K9").Value = 1.45

$ws.Range("B10").Value = 1.27
$ws.Range("E10").Value = 0.95
$ws.Range("H10").Value = 2.16
$ws.Range("K10").Value = 1.51

$ws.Range("B11").Value = 1.32
$ws.Range("E11").Value = 1.01
$ws.Range("H11").Value = 2.25
$ws.Range("K11").Value = 1.58

$ws.Range("B12").Value = 1.38
$ws.Range("E12").Value = 1.07
$ws.Range("H12").Value = 2.34
$ws.Range("K12").Value = 1.65

$ws.Range("B13").Value = 1.45
$ws.Range("E13").Value = 1.13
$ws.Range("H13").Value = 2.43
$ws.Range("K13").Value = 1.73

$ws.Range("B14").Value = 1.53
$ws.Range("E14").Value = 1.2
$ws.Range("H14").Value = 2.54
$ws.Range("K14").Value = 1.8

$ws.Range("B15").Value = 1.62
$ws.Range("E15").Value = 1.28
$ws.Range("H15").Value = 2.65
$ws.Range("K15").Value = 1.89

$ws.Range("B16").Value = 1.72
$ws.Range("E16").Value = 1.37
$ws.Range("H16").Value = 2.76
$ws.Range("K16").Value = 1.97

$ws.Range("B17").Value = 1.84
$ws.Range("E17").Value = 1.46
$ws.Range("H17").Value = 2.89
$ws.Range("K17").Value = 2.06

$ws.Range("B18").Value = 1.97
$ws.Range("E18").Value = 1.57
$ws.Range("H18").Value = 3.02
$ws.Range("K18").Value = 2.16

$ws.Range("B19").Value = 2.13
$ws.Range("E19").Value = 1.68
$ws.Range("H19").Value = 3.16
$ws.Range("K19").Value = 2.26

$ws.Range("B20").Value = 2.31
$ws.Range("E20").Value = 1.8
$ws.Range("H20").Value = 3.31
$ws.Range("K20").Value = 2.36

$ws.Range("B21").Value = 2.52
$ws.Range("E21").Value = 1.93
$ws.Range("H21").Value = 3.47
$ws.Range("K21").Value = 2.47

$ws.Range("B22").Value = 2.75
$ws.Range("E22").Value = 2.08
$ws.Range("H22").Value = 3.64
$ws.Range("K22").Value = 2.57

$ws.Range("B23").Value = 3.01
$ws.Range("E23").Value = 2.23
$ws.Range("H23").Value = 3.82
$ws.Range("K23").Value = 2.68

$ws.Range("B24").Value = 3.3
$ws.Range("E24").Value = 2.4
$ws.Range("H24").Value = 4.01
$ws.Range("K24").Value = 2.79

$ws.Range("B25").Value = 3.62
$ws.Range("E25").Value = 2.58
$ws.Range("H25").Value = 4.2
$ws.Range("K25").Value = 2.91

$ws.Range("B26").Value = 3.98
$ws.Range("E26").Value = 2.79
$ws.Range("H26").Value = 4.4
$ws.Range("K26").Value = 3.02

$ws.Range("B27").Value = 4.37
$ws.Range("E27").Value = 3.01
$ws.Range("H27").Value = 4.6
$ws.Range("K27").Value = 3.14

$ws.Range("B28").Value = 4.8
$ws.Range("E28").Value = 3.25
$ws.Range("H28").Value = 4.8
$ws.Range("K28").Value = 3.25

$ws.Range("B29").Value = 5.28
$ws.Range("E29").Value = 3.52
$ws.Range("H29").Value = 5
$ws.Range("K29").Value = 3.37

$ws.Range("B30").Value = 5.8
$ws.Range("E30").Value = 3.82
$ws.Range("H30").Value = 5.2
$ws.Range("K30").Value = 3.49

$ws.Range("B31").Value = 6.37
$ws.Range("E31").Value = 4.15
$ws.Range("H31").Value = 5.39
$ws.Range("K31").Value = 3.6

$ws.Range("B32").Value = 7
$ws.Range("E32").Value = 4.52
$ws.Range("H32").Value = 5.56
$ws.Range("K32").Value = 3.71

$ws.Range("B33").Value = 7.68
$ws.Range("E33").Value = 4.92
$ws.Range("H33").Value = 5.72
$ws.Range("K33").Value = 3.81

$ws.Range("B34").Value = 8.43
$ws.Range("E34").Value = 5.37
$ws.Range("H34").Value = 5.87
$ws.Range("K34").Value = 3.91

$ws.Range("B35").Value = 9.24
$ws.Range("E35").Value = 5.87
$ws.Range("H35").Value = 6.01
$ws.Range("K35").Value = 4

$ws.Range("B36").Value = 10.12
$ws.Range("E36").Value = 6.41
$ws.Range("H36").Value = 6.14
$ws.Range("K36").Value = 4.1

$ws.Range("B37").Value = 11.08
$ws.Range("E37").Value = 7.02
$ws.Range("H37").Value = 6.26
$ws.Range("K37").Value = 4.19

$ws.Range("B38").Value = 12.12
$ws.Range("E38").Value = 7.68
$ws.Range("H38").Value = 6.37
$ws.Range("K38").Value = 4.29

